$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update "Marking" row: Right count 3 -> 5
$ws.Range("B11").Value = 5

# Update "Total" row: Right count 57 -> 95
$ws.Range("B12").Value = 95

# Update "Total" row Max column: "54/84" -> "95/140"
$ws.Range("E12").Value = "95/140"
